$d = $word.ActiveDocument

# 1. Remove the leading "01" run from the "Week Ending" date paragraph,
#    turning "01 September 2023" into " September 2023".
$d.Content.Find.Execute("01", $true, $true, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 2. Re-apply bullet list formatting across the three "Accomplishments" bullet
#    paragraphs (Continued with milestone 2, / Completed the simulation for
#    moving 1 meter.../ Completed the simulation for turning 90, 180, and 270
#    degrees.) so Word mints a fresh list definition for them (numId 1 -> 4).
$pStart = $d.Paragraphs(23)
$pEnd = $d.Paragraphs(25)
$bulletRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$bulletRange.ListFormat.ApplyBulletDefault()

# 3. Merge "Continued with" + " milestone 2, " into a single run.
$d.Content.Find.Execute("Continued with milestone 2, ", $true, $false, $false, $false, $false, $true, 1, $false, "Continued with milestone 2, ", 2) | Out-Null

# 4. Drop the trailing space on the "turning 90, 180, and 270 degrees" bullet.
$d.Content.Find.Execute("Completed the simulation for turning 90, 180, and 270 degrees. ", $true, $false, $false, $false, $false, $true, 1, $false, "Completed the simulation for turning 90, 180, and 270 degrees.", 2) | Out-Null

# 5. Drop the trailing space on "Complete the rest of the simulation files."
$d.Content.Find.Execute("Complete the rest of the simulation files. ", $true, $false, $false, $false, $false, $true, 1, $false, "Complete the rest of the simulation files.", 2) | Out-Null

# 6. Merge "Get a start on the " + "intersection algorithm and distance measuring" + ". "
#    into a single run with no trailing space.
$d.Content.Find.Execute("Get a start on the intersection algorithm and distance measuring. ", $true, $false, $false, $false, $false, $true, 1, $false, "Get a start on the intersection algorithm and distance measuring.", 2) | Out-Null

Write-Output "Week Ending paragraph: [$($d.Paragraphs(4).Range.Text)]"
Write-Output "Bullet 1: [$($d.Paragraphs(23).Range.Text)]"
Write-Output "Bullet 2: [$($d.Paragraphs(24).Range.Text)]"
Write-Output "Bullet 3: [$($d.Paragraphs(25).Range.Text)]"
Write-Output "Upcoming 1: [$($d.Paragraphs(32).Range.Text)]"
Write-Output "Upcoming 2: [$($d.Paragraphs(33).Range.Text)]"
